$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quellen DigCommTool")

# Two new source entries (icons used in the UI), each described by a
# Beschreibung / URL / Zugriffsdatum triple, matching the existing table
# layout in columns A:C.

# URLs (column B) for the new rows 4 and 5.
$ws.Range("B4").Value = "https://www.clipartmax.com/png/small/330-3306271_wall-icon-clipart-wall-brick-building-wall-icon.png"
$ws.Range("B5").Value = "https://www.freeiconspng.com/img/33857"

# Descriptions (column A).
$ws.Range("A4").Value = "Icon"
$ws.Range("A5").Value = "Icon"

# Access dates (column C) - entered as literal text ("08.01.21") rather
# than a parsed date serial, so write the text via a formula and collapse
# it to a plain value (keeps the cell's existing date-formatted style
# untouched instead of minting a new number-format style).
$ws.Range("C4").Formula = "=""08.01.21"""
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("C5").Formula = "=""08.01.21"""
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Update selection to G16
$ws.Range("G16").Select()
